# Add the GitHub repo URL (as a hyperlink) to the title slide's subtitle,
# fix the "dawoud" -> "Dawoud" capitalization, and give the subtitle shape
# an explicit position/size with "shrink text on overflow" autofit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)   # "Subtitle 2" placeholder (subTitle idx=1)

# --- Explicit position / size for the subtitle placeholder ---------------
# Shape.Left/Top/Width/Height are in points; the host stores them as
# single-precision floats before converting to EMU (1 pt = 12700 EMU), so
# the literals below are chosen to land exactly on the target EMU values
# (2417780, 3531204, 8637072, 2085825) once they round-trip through that
# float32 conversion.
$shape.Left   = 190.3763809527559
$shape.Top    = 278.0475616551181
$shape.Width  = 680.0844094488189
$shape.Height = 164.23818897637796

# --- Shrink text on overflow (<a:bodyPr><a:normAutofit/></a:bodyPr>) -----
$shape.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape

# --- Fix "dawoud" -> "Dawoud" and clear the flagged-as-misspelled marker -
# Deleting the run and re-inserting fresh text (rather than editing the
# existing run's .Text in place) produces a clean <a:rPr> with no err="1".
$tr = $shape.TextFrame.TextRange
$fullText = $tr.Text
$idx = $fullText.IndexOf("dawoud")
$oldRun = $tr.Characters($idx + 1, 6)
$oldRun.Delete() | Out-Null
$shape.TextFrame.TextRange.InsertAfter("Dawoud") | Out-Null

# --- Add a blank paragraph, then a paragraph with the GitHub URL ---------
$cr = [char]13
$shape.TextFrame.TextRange.InsertAfter($cr) | Out-Null
$shape.TextFrame.TextRange.InsertAfter($cr + "https://github.com/theRadFad/HeartMonitor") | Out-Null

# --- Turn the newly-added URL text into a hyperlink -----------------------
$urlText = "https://github.com/theRadFad/HeartMonitor"
$full = $shape.TextFrame.TextRange
$urlStart = $full.Length - $urlText.Length + 1
$urlRange = $full.Characters($urlStart, $urlText.Length)
$urlRange.ActionSettings(1).Hyperlink.Address = $urlText
